# "switch to system user prompts"
# ContractAIForm (sheet19.xml): insert two new columns - "systemPrompt" and
# "userPrompt" - right before the existing "form" column (which, along with
# every column after it, shifts two places to the right: F:K -> H:M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ContractAIForm")

# Insert two blank columns at F:G - everything from the old column F
# ("form") onward shifts right by two columns.
$ws.Range("F1:G1").EntireColumn.Insert()

# Row 1 - column headers
$ws.Range("F1").Value = "systemPrompt"
$ws.Range("G1").Value = "userPrompt"

# Row 2 - sqlType
$ws.Range("F2").Value = "VARCHAR(255)"
$ws.Range("G2").Value = "VARCHAR(255)"

# Row 3 - isRequired (blank for these two new columns)
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""

# Row 4 - refernces (blank for these two new columns)
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
